$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 14, pushing the existing "Total" row (and its formulas) down to row 15.
$ws.Rows("14:14").Insert()

# Populate the new row 14 with a new timesheet entry (same From/To/Rate as the entry above it).
$ws.Range("A14").Value = 45278
$ws.Range("B14").Value = 0.583333333333333
$ws.Range("C14").Value = 0.833333333333333
$ws.Range("D14").Formula = "=(C14<B14)+C14-B14"
$ws.Range("E14").Value = 10
$ws.Range("F14").Formula = "=(D14*24)*E14"

# Copy the formatting from the row above (row 13) into the new row 14 (applied last, so the
# auto-number-formatting that Excel performs when entering formulas doesn't override it), so
# it keeps the same per-column styles (date / time / hours / rate columns) as the rest of the table.
$ws.Range("A13:F13").Copy()
$ws.Range("A14:F14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fix up the Total row (now row 15) so it sums through the new row 14.
$ws.Range("D15").Formula = "=SUM(D2:D14)"
$ws.Range("F15").Formula = "=SUM(F2:F14)"

$ws.Range("F16").Select()
